$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.169.50'
$ws.Range('D2').Style = $__style
$ws.Range('E2').Value = '  -1.07%  '
$__style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.660.17'
$ws.Range('D3').Style = $__style
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  +0.29%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.48'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -1.15%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5207'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E7').Value = '  +0.27%  '
$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2639'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  -2.58%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06282'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  -1.96%  '
$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.80'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  -4.87%  '
$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07727'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$__style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.431'
$ws.Range('D12').Style = $__style
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.642.14'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  -2.12%  '
$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5439'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  -2.73%  '
$__style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0₅8171'
$ws.Range('D16').Style = $__style
$ws.Range('E16').Value = '  -1.84%  '
$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.52'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  -1.78%  '
$__style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '26.199.70'
$ws.Range('D18').Style = $__style
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('E19').Value = '  +0.39%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.632'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  -3.46%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '192.17'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  -0.62%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.08'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -2.14%  '
$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.062'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  -4.14%  '
$ws.Range('E24').Value = '  +0.45%  '
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '139.95'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -0.54%  '
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1230'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  -3.86%  '
$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.181'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  -3.09%  '
$ws.Range('E28').Value = '  -1.54%  '
$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.407'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  -2.68%  '
$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05995'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  -4.58%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.271'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  -1.03%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.553'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -1.50%  '
$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.257'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  -5.84%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.611'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -5.18%  '
$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9677'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  -4.34%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.769'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  -0.49%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5675'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  -7.74%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.012'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01597'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  -2.09%  '
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8559'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  +0.34%  '
$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.016.18'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  -7.39%  '
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.56'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  -0.10%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.800.81'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  -1.20%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '56.83'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₈107'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  -5.45%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.008'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  +0.88%  '
$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.006'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -2.29%  '
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05173'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  -0.69%  '
$__style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.456'
$ws.Range('D51').Style = $__style
$ws.Range('E51').Value = '  -1.20%  '
